$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")
$ws.Range("D1").Value = "Remaining time"
$ws.Range("H3").Value = 43
$ws.Range("H4").Value = -28
$ws.Range("H5").Value = 21
$ws.Range("H6").Value = -8
